# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in columns AD, AE, AF (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) used by A1:AC1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the team's record for every player row (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
